$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "37.407.49"
$ws.Cells.Item(3, 4).Value = "2.079.96"
$ws.Cells.Item(3, 5).Value = "  -2.17%  "
$ws.Cells.Item(4, 5).Value = "  +0.18%  "
Set-TextValue 5 4 "251.25"
$ws.Cells.Item(5, 5).Value = "  +0.29%  "
$ws.Cells.Item(7, 5).Value = "  -0.07%  "
Set-TextValue 8 4 "56.86"
$ws.Cells.Item(8, 5).Value = "  +25.51%  "
Set-TextValue 9 4 "62.11"
$ws.Cells.Item(9, 5).Value = "  +1.42%  "
Set-TextValue 10 4 "0.388"
$ws.Cells.Item(10, 5).Value = "  +5.23%  "
Set-TextValue 11 4 "0.0758"
$ws.Cells.Item(11, 5).Value = "  +3.17%  "
Set-TextValue 12 4 "0.107"
$ws.Cells.Item(12, 5).Value = "  +7.21%  "
Set-TextValue 13 4 "15.62"
$ws.Cells.Item(13, 5).Value = "  +6.53%  "
$ws.Cells.Item(14, 4).Value = "2.382.52"
$ws.Cells.Item(14, 5).Value = "  -1.73%  "
Set-TextValue 15 4 "0.846"
$ws.Cells.Item(15, 5).Value = "  -0.63%  "
$ws.Cells.Item(16, 5).Value = "  +4.45%  "
$ws.Cells.Item(17, 4).Value = "2.079.16"
$ws.Cells.Item(17, 5).Value = "  -2.16%  "
$ws.Cells.Item(18, 4).Value = "37.332.77"
$ws.Cells.Item(18, 5).Value = "  +1.36%  "
Set-TextValue 19 4 "73.29"
$ws.Cells.Item(19, 5).Value = "  -0.26%  "
Set-TextValue 20 4 "14.80"
$ws.Cells.Item(20, 5).Value = "  +13.90%  "
$ws.Cells.Item(21, 4).Value = "0.0₃0852"
$ws.Cells.Item(21, 5).Value = "  +3.57%  "
Set-TextValue 22 4 "241.25"
$ws.Cells.Item(22, 5).Value = "  -0.02%  "
$ws.Cells.Item(23, 5).Value = "  +3.25%  "
$ws.Cells.Item(24, 5).Value = "  +0.13%  "
$ws.Cells.Item(25, 5).Value = "  +0.03%  "
Set-TextValue 26 4 "171.86"
Set-TextValue 27 4 "9.23"
$ws.Cells.Item(27, 5).Value = "  +2.09%  "
Set-TextValue 28 4 "21.05"
$ws.Cells.Item(28, 5).Value = "  +1.65%  "
$ws.Cells.Item(29, 5).Value = "  +0.88%  "
$ws.Cells.Item(30, 5).Value = "  +0.68%  "
Set-TextValue 31 4 "23.43"
$ws.Cells.Item(31, 5).Value = "  +5.79%  "
$ws.Cells.Item(32, 5).Value = "  +22.42%  "
Set-TextValue 33 4 "4.59"
$ws.Cells.Item(33, 5).Value = "  +2.82%  "
Set-TextValue 34 4 "0.0631"
$ws.Cells.Item(34, 5).Value = "  +5.38%  "
Set-TextValue 35 4 "4.37"
$ws.Cells.Item(35, 5).Value = "  +7.05%  "
Set-TextValue 36 4 "0.0910"
$ws.Cells.Item(36, 5).Value = "  -0.08%  "
$ws.Cells.Item(37, 5).Value = "  +0.05%  "
$ws.Cells.Item(38, 2).Value = "LidoDAOToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue 38 4 "2.30"
$ws.Cells.Item(38, 5).Value = "  -1.27%  "
$ws.Cells.Item(39, 2).Value = "WEMIXToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue 39 4 "1.86"
$ws.Cells.Item(39, 5).Value = "  -1.30%  "
$ws.Cells.Item(40, 5).Value = "  -0.30%  "
$ws.Cells.Item(41, 2).Value = "Cronos"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue 41 4 "0.101"
$ws.Cells.Item(41, 5).Value = "  +19.75%  "
$ws.Cells.Item(42, 2).Value = "VeChain"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue 42 4 "0.0229"
$ws.Cells.Item(42, 5).Value = "  +4.38%  "
$ws.Cells.Item(43, 5).Value = "  +8.76%  "
Set-TextValue 44 4 "1.18"
$ws.Cells.Item(44, 5).Value = "  -2.45%  "
Set-TextValue 45 4 "100.25"
$ws.Cells.Item(45, 5).Value = "  -0.31%  "
Set-TextValue 46 4 "2.81"
$ws.Cells.Item(46, 5).Value = "  -0.21%  "
Set-TextValue 47 4 "4.17"
$ws.Cells.Item(47, 5).Value = "  +89.41%  "
$ws.Cells.Item(48, 4).Value = "1.327.85"
$ws.Cells.Item(48, 5).Value = "  -2.77%  "
Set-TextValue 49 4 "2.41"
$ws.Cells.Item(49, 5).Value = "  +5.38%  "
Set-TextValue 50 4 "2.93"
$ws.Cells.Item(50, 5).Value = "  +3.19%  "
Set-TextValue 51 4 "7.05"
$ws.Cells.Item(51, 5).Value = "  +8.09%  "
